$wb = $excel.ActiveWorkbook

# --- "My Side" sheet: drop the trailing GST-portal-matched columns (J:Q) ---
$mySide = $wb.Worksheets.Item("My Side")
$mySide.Range("J1:Q23").Delete()

# --- "GST portal" sheet: drop the not-matching columns and shift the
#     remaining ones (K:Q) left so they sit right after the unique id in A ---
$gst = $wb.Worksheets.Item("GST portal")
$gst.Range("J1").EntireColumn.Delete()
$gst.Range("I1").EntireColumn.Delete()
$gst.Range("H1").EntireColumn.Delete()
$gst.Range("G1").EntireColumn.Delete()
$gst.Range("F1").EntireColumn.Delete()
$gst.Range("E1").EntireColumn.Delete()
$gst.Range("D1").EntireColumn.Delete()
$gst.Range("B1").EntireColumn.Delete()
